# Update harmonized emissions values (regression test baseline update for
# actual base_year) on the "data" worksheet, rows 2-7, columns G-O.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (World, sector1, SSP scenario A)
$ws.Range("G2").Value = 32.01690590111642
$ws.Range("H2").Value = 32.98149920255183
$ws.Range("I2").Value = 33.65986177565124
$ws.Range("J2").Value = 34.05199362041467
$ws.Range("K2").Value = 34.15789473684211
$ws.Range("L2").Value = 35.69473684210526
$ws.Range("M2").Value = 37.14736842105263
$ws.Range("N2").Value = 38.51578947368422
$ws.Range("O2").Value = 39.8

# Row 3 (World, sector2, SSP scenario A)
$ws.Range("G3").Value = 15.16842105263158
$ws.Range("H3").Value = 16.04210526315789
$ws.Range("I3").Value = 16.83157894736842
$ws.Range("J3").Value = 17.53684210526316
$ws.Range("K3").Value = 18.15789473684211
$ws.Range("L3").Value = 18.69473684210526
$ws.Range("M3").Value = 19.14736842105263
$ws.Range("N3").Value = 19.51578947368421
$ws.Range("O3").Value = 19.8

# Row 4 (World, total, SSP scenario A)
$ws.Range("G4").Value = 16.84848484848485
$ws.Range("H4").Value = 16.93939393939394
$ws.Range("I4").Value = 16.82828282828283
$ws.Range("J4").Value = 16.51515151515152

# Row 5 (World, sector1, SSP scenario B)
$ws.Range("G5").Value = 32.01690590111642
$ws.Range("H5").Value = 32.98149920255183
$ws.Range("I5").Value = 33.65986177565124
$ws.Range("J5").Value = 34.05199362041467
$ws.Range("K5").Value = 34.15789473684211
$ws.Range("L5").Value = 35.69473684210526
$ws.Range("M5").Value = 37.14736842105263
$ws.Range("N5").Value = 38.51578947368422
$ws.Range("O5").Value = 39.8

# Row 6 (World, sector2, SSP scenario B)
$ws.Range("G6").Value = 15.16842105263158
$ws.Range("H6").Value = 16.04210526315789
$ws.Range("I6").Value = 16.83157894736842
$ws.Range("J6").Value = 17.53684210526316
$ws.Range("K6").Value = 18.15789473684211
$ws.Range("L6").Value = 18.69473684210526
$ws.Range("M6").Value = 19.14736842105263
$ws.Range("N6").Value = 19.51578947368421
$ws.Range("O6").Value = 19.8

# Row 7 (World, total, SSP scenario B)
$ws.Range("G7").Value = 16.84848484848485
$ws.Range("H7").Value = 16.93939393939394
$ws.Range("I7").Value = 16.82828282828283
$ws.Range("J7").Value = 16.51515151515152
